# Natmi following Dr Hou advice
# Update Ligand-expressing cells, Receptor-expressing cells and all derived
# expression / specificity metrics for rows 2-4 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 106.3055973333333
$ws.Range("H2").Value = 318.916792
$ws.Range("I2").Value = 0.2547398208373942
$ws.Range("J2").Value = 0.2547398208373942
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.321929333333333
$ws.Range("N2").Value = 21.965788
$ws.Range("Q2").Value = 778.3620714124551
$ws.Range("R2").Value = 7005.258642712096
$ws.Range("S2").Value = 0.2547398208373942
$ws.Range("T2").Value = 0.2547398208373942

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 296.1091513333333
$ws.Range("H3").Value = 888.327454
$ws.Range("I3").Value = 0.7095655736964096
$ws.Range("J3").Value = 0.7095655736964096
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.321929333333333
$ws.Range("N3").Value = 21.965788
$ws.Range("Q3").Value = 2168.090281015972
$ws.Range("R3").Value = 19512.81252914375
$ws.Range("S3").Value = 0.7095655736964096
$ws.Range("T3").Value = 0.7095655736964096

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.895733
$ws.Range("H4").Value = 44.68719899999999
$ws.Range("I4").Value = 0.03569460546619627
$ws.Range("J4").Value = 0.03569460546619627
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.321929333333333
$ws.Range("N4").Value = 21.965788
$ws.Range("Q4").Value = 109.0655043942013
$ws.Range("R4").Value = 981.5895395478118
$ws.Range("S4").Value = 0.03569460546619627
$ws.Range("T4").Value = 0.03569460546619627
